$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 46056
$ws.Cells.Item(2, 3).Value = 0.190988057015924
$ws.Cells.Item(3, 1).Value = 46057
$ws.Cells.Item(3, 3).Value = 0.1877706628093372
$ws.Cells.Item(4, 1).Value = 46058
$ws.Cells.Item(4, 3).Value = 0.1895076263294159
$ws.Cells.Item(5, 1).Value = 46059
$ws.Cells.Item(5, 3).Value = 0.1844419478869654
$ws.Cells.Item(6, 1).Value = 46062
$ws.Cells.Item(6, 3).Value = 0.187153842966375
$ws.Cells.Item(7, 1).Value = 46063
$ws.Cells.Item(7, 3).Value = 0.1858241617800606
$ws.Cells.Item(8, 1).Value = 46064
$ws.Cells.Item(8, 3).Value = 0.1864312574827131
$ws.Cells.Item(9, 1).Value = 46065
$ws.Cells.Item(9, 3).Value = 0.1858347525854054
$ws.Cells.Item(10, 1).Value = 46066
$ws.Cells.Item(10, 3).Value = 0.1845732316294208
$ws.Cells.Item(11, 1).Value = 46069
$ws.Cells.Item(11, 3).Value = 0.1844089816523485
$ws.Cells.Item(12, 1).Value = 46070
$ws.Cells.Item(12, 3).Value = 0.1817907350685866
$ws.Cells.Item(13, 1).Value = 46071
$ws.Cells.Item(13, 3).Value = 0.1828346324605994
$ws.Cells.Item(14, 1).Value = 46072
$ws.Cells.Item(14, 3).Value = 0.1836150703322842
$ws.Cells.Item(15, 1).Value = 46073
$ws.Cells.Item(15, 3).Value = 0.1791652003320497
$ws.Cells.Item(16, 1).Value = 46076
$ws.Cells.Item(16, 3).Value = 0.1805293931913482
$ws.Cells.Item(17, 1).Value = 46077
$ws.Cells.Item(17, 3).Value = 0.1797816846968144
$ws.Cells.Item(18, 1).Value = 46078
$ws.Cells.Item(18, 3).Value = 0.1774258281060843
$ws.Cells.Item(19, 1).Value = 46079
$ws.Cells.Item(19, 3).Value = 0.1806075806058409
$ws.Cells.Item(20, 1).Value = 46080
$ws.Cells.Item(20, 3).Value = 0.182799531969842
$ws.Cells.Item(21, 1).Value = 46083
$ws.Cells.Item(21, 3).Value = 0.1810470457841293
$ws.Cells.Item(22, 1).Value = 46084
$ws.Cells.Item(22, 3).Value = 0.1795510363419213
$ws.Cells.Item(23, 1).Value = 46085
$ws.Cells.Item(23, 3).Value = 0.183205352023283
$ws.Cells.Item(24, 1).Value = 46086
$ws.Cells.Item(24, 3).Value = 0.1851226691744886
$ws.Cells.Item(25, 1).Value = 46087
$ws.Cells.Item(25, 3).Value = 0.1871951897285543
$ws.Cells.Item(26, 1).Value = 46090
$ws.Cells.Item(26, 3).Value = 0.1905490388915473
$ws.Cells.Item(27, 1).Value = 46091
$ws.Cells.Item(27, 3).Value = 0.1915829261951968
$ws.Cells.Item(28, 1).Value = 46092
$ws.Cells.Item(28, 3).Value = 0.1922165870240487
$ws.Cells.Item(29, 1).Value = 46093
$ws.Cells.Item(29, 3).Value = 0.1927710666461157
$ws.Cells.Item(30, 1).Value = 46094
$ws.Cells.Item(30, 3).Value = 0.1928736015362905
$ws.Cells.Item(31, 1).Value = 46097
$ws.Cells.Item(31, 3).Value = 0.1929361193580793
